$p = $ppt.ActivePresentation

# Remove the second slide (SlideID 256) - the "x = 1 / print(x)" variable
# demo slide with the grey/red rectangle boxes. PowerPoint renumbers the
# remaining slide/master relationship ids automatically on save.
$p.Slides.Item(2).Delete()
